# Update "best parameters" summary tables + underlying stats after
# regenerating the graphs/graph code (per commit message).
#
# Net effect (verified against the canonical OOXML diff): only numeric
# cell values change on four sheets. The workbook's <mergeCells> lists
# are re-serialised in a different order in the target file, but the
# *set* of merged ranges on every sheet is unchanged - that's an
# artifact of whatever tool re-saved the file, not a deliberate edit,
# so there is nothing to replicate there.

$wb = $excel.ActiveWorkbook

# --- Sheet: Manual calibration (sum) ---------------------------------
$ws1 = $wb.Worksheets.Item("Manual calibration (sum)")
$ws1.Range("B4").Value = 150.0
$ws1.Range("C4").Value = 15.0
$ws1.Range("F4").Value = 93.333333
$ws1.Range("G4").Value = -0.383575
$ws1.Range("H4").Value = -1.518704
$ws1.Range("I4").Value = -2.589177

# --- Sheet: Manual calibration (indiv) -------------------------------
$ws2 = $wb.Worksheets.Item("Manual calibration (indiv)")
$ws2.Range("B4").Value = 150.0
$ws2.Range("C4").Value = 15.0
$ws2.Range("F4").Value = 93.333333
$ws2.Range("G4").Value = -0.383575
$ws2.Range("H4").Value = -1.518704
$ws2.Range("I4").Value = -2.589177

# --- Sheet: BOTorch (sum) --------------------------------------------
$ws3 = $wb.Worksheets.Item("BOTorch (sum)")
$ws3.Range("K6").Value = 16.438385
$ws3.Range("L6").Value = 7.963184
$ws3.Range("O6").Value = 206.411129
$ws3.Range("P6").Value = -0.573375
$ws3.Range("Q6").Value = -2.046668
$ws3.Range("R6").Value = -3.174586
$ws3.Range("S6").Value = -1.931543

# --- Sheet: BOTorch (indiv) -------------------------------------------
$ws4 = $wb.Worksheets.Item("BOTorch (indiv)")
$ws4.Range("K12").Value = 16.438385
$ws4.Range("L12").Value = 7.963184
$ws4.Range("O12").Value = 206.411129
$ws4.Range("P12").Value = -0.573375
$ws4.Range("Q12").Value = -2.046668
$ws4.Range("R12").Value = -3.174586
$ws4.Range("S12").Value = -1.931543
